$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the cryptos list (Price + Volume(1h) columns) with the latest
# scraped figures. A handful of "Price" cells are plain decimals
# (e.g. "1.00", "198.78") that Excel would otherwise auto-convert to
# numbers on assignment, dropping the trailing zero / losing the exact
# text, so those use a leading apostrophe to force literal text, just
# like a user typing '1.00 directly into the cell.
$ws.Range("D2").Value = "76.512.19"
$ws.Range("E2").Value = "  +1.01%  "
$ws.Range("D3").Value = "2.963.71"
$ws.Range("E3").Value = "  +3.18%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'198.78"
$ws.Range("E5").Value = "  +1.94%  "
$ws.Range("D6").Value = "'596.41"
$ws.Range("E6").Value = "  -0.18%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").Value = "'0.205"
$ws.Range("E9").Value = "  +7.24%  "
$ws.Range("D10").Value = "2.964.35"
$ws.Range("E10").Value = "  +3.14%  "
$ws.Range("D11").Value = "'0.443"
$ws.Range("E11").Value = "  +10.84%  "
$ws.Range("E12").Value = "  +0.53%  "
$ws.Range("B13").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C13").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D13").Value = "3.509.61"
$ws.Range("E13").Value = "  +2.44%  "
$ws.Range("B14").Value = "Toncoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D14").Value = "'4.92"
$ws.Range("E14").Value = "  +0.15%  "
$ws.Range("D15").Value = "'28.58"
$ws.Range("E15").Value = "  +4.68%  "
$ws.Range("D16").Value = "76.449.12"
$ws.Range("E16").Value = "  +0.95%  "
$ws.Range("D17").Value = "'0.0000191"
$ws.Range("E17").Value = "  +1.03%  "
$ws.Range("D18").Value = "2.956.07"
$ws.Range("E18").Value = "  +1.97%  "
$ws.Range("D19").Value = "'13.66"
$ws.Range("E19").Value = "  +8.79%  "
$ws.Range("E20").Value = "  -1.95%  "
$ws.Range("D21").Value = "'378.82"
$ws.Range("E21").Value = "  -0.54%  "
$ws.Range("D22").Value = "'2.29"
$ws.Range("E22").Value = "  +0.14%  "
$ws.Range("E23").Value = "  +4.79%  "
$ws.Range("D24").Value = "'72.36"
$ws.Range("E24").Value = "  +0.99%  "
$ws.Range("B25").Value = "Dai"
$ws.Range("C25").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D25").Value = "'1.00"
$ws.Range("E25").Value = "  +0.26%  "
$ws.Range("B26").Value = "WrappedeETH"
$ws.Range("C26").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D26").Value = "3.093.38"
$ws.Range("E26").Value = "  +1.67%  "
$ws.Range("D27").Value = "'4.31"
$ws.Range("E27").Value = "  +2.14%  "
$ws.Range("D28").Value = "'9.77"
$ws.Range("D29").Value = "'0.0000108"
$ws.Range("E29").Value = "  +1.11%  "
$ws.Range("E30").Value = "  +0.16%  "
$ws.Range("D31").Value = "'8.60"
$ws.Range("E31").Value = "  +10.57%  "
$ws.Range("E32").Value = "  -1.16%  "
$ws.Range("D33").Value = "'497.12"
$ws.Range("E33").Value = "  -1.93%  "
$ws.Range("D34").Value = "'1.82"
$ws.Range("E34").Value = "  +0.44%  "
$ws.Range("E35").Value = "  +0.04%  "
$ws.Range("B36").Value = "Monero"
$ws.Range("C36").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D36").Value = "'164.75"
$ws.Range("E36").Value = "  +0.10%  "
$ws.Range("B37").Value = "EthereumClassic"
$ws.Range("C37").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D37").Value = "'20.33"
$ws.Range("E37").Value = "  +1.01%  "
$ws.Range("E38").Value = "  +14.79%  "
$ws.Range("E39").Value = "  +19.23%  "
$ws.Range("D40").Value = "'19.97"
$ws.Range("E40").Value = "  +1.55%  "
$ws.Range("E41").Value = "  -2.20%  "
$ws.Range("E42").Value = "  +0.04%  "
$ws.Range("D43").Value = "'180.78"
$ws.Range("E43").Value = "  -1.30%  "
$ws.Range("E44").Value = "  -0.98%  "
$ws.Range("E45").Value = "  -1.26%  "
$ws.Range("D46").Value = "'40.01"
$ws.Range("E46").Value = "  -0.47%  "
$ws.Range("E47").Value = "  -1.89%  "
$ws.Range("E48").Value = "  +2.97%  "
$ws.Range("E49").Value = "  +4.11%  "
$ws.Range("D50").Value = "'2.33"
$ws.Range("E50").Value = "  -0.89%  "
$ws.Range("D51").Value = "'0.669"
$ws.Range("E51").Value = "  +0.42%  "
